# Updated symbol list on Thu Feb  2 21:34:10 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns with the latest scraped
# values. These columns store plain text (not numbers/percent-numbers), so
# we force a Text number format before writing the strings to avoid Excel
# auto-converting them into numeric/percentage values (which would corrupt
# values like "8.780" -> 8.78 or "3.76%" -> 0.0376).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "328.83";    "E2"  = "3.76%"
    "D3"  = "40.37";     "E3"  = "6.96%"
    "D4"  = "5.827";     "E4"  = "12.93%"
    "D5"  = "0.08066";   "E5"  = "1.20%"
                         "E6"  = "3.04%"
    "D7"  = "8.780";     "E7"  = "3.28%"
    "D8"  = "1.960";     "E8"  = "1.22%"
    "D9"  = "2.941";     "E9"  = "-1.19%"
    "D10" = "0.9444";    "E10" = "0.46%"
                         "E11" = "0.92%"
    "D12" = "0.1984";    "E12" = "2.83%"
    "D13" = "8.929";     "E13" = "36.99%"
    "D14" = "0.09205";   "E14" = "2.01%"
    "D15" = "0.03510";   "E15" = "3.42%"
    "D16" = "0.09613";   "E16" = "1.01%"
    "D17" = "0.001312";  "E17" = "-3.91%"
    "D18" = "0.006045";  "E18" = "0.30%"
    "D19" = "3.370";     "E19" = "-1.20%"
    "D20" = "0.3564";    "E20" = "1.46%"
    "D21" = "0.1432";    "E21" = "9.79%"
    "D22" = "0.2415";    "E22" = "4.99%"
    "D23" = "0.04427";   "E23" = "1.86%"
    "D24" = "0.001260";  "E24" = "5.28%"
    "D25" = "0.004378";  "E25" = "-0.90%"
    "D26" = "0.0001143"; "E26" = "-13.56%"
                         "E27" = "0.67%"
    "D39" = "0.02453";   "E39" = "3.57%"
    "D40" = "0.05313";   "E40" = "2.72%"
    "D41" = "0.007433";  "E41" = "0.35%"
    "D42" = "0.1420";    "E42" = "1.92%"
    "D43" = "0.008731";  "E43" = "2.90%"
    "D44" = "0.002107";  "E44" = "5.81%"
    "D45" = "0.01085";   "E45" = "23.97%"
    "D46" = "0.00006898";"E46" = "7.85%"
                         "E47" = "0.87%"
    "D48" = "0.003176";  "E48" = "11.47%"
                         "E49" = "1.50%"
                         "E50" = "0.87%"
    "D51" = "0.0002006"; "E51" = "0.87%"
}

# Force each target cell to the Text number format before writing so Excel
# keeps the strings verbatim instead of re-interpreting them as numbers or
# percentages (which would corrupt values like "8.780" -> 8.78 or
# "3.76%" -> 0.0376). A Union-range NumberFormat assignment only affects the
# first area under this host, so format + write cell-by-cell instead.
foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
